$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "mem2Reg" instruction row (row 9) is being dropped from the truth
# table entirely.
$ws.Rows(9).Delete()

# Add the new "comp" (comparator) instruction as column O: a header cell
# plus the truth-table bit for each control-signal row (branch, regw,
# regdes, alusrc, memR, memW, baw).
$compValues = @("comp", 0, 1, 1, 0, 0, 0, 0)

for ($i = 0; $i -lt $compValues.Count; $i++) {
    $row = $i + 1
    $cell = $ws.Range("O$row")
    $cell.Value = $compValues[$i]
    # Match the centered alignment used by the rest of the table.
    $cell.HorizontalAlignment = -4108
}

# Match the selection left behind in the saved file.
$ws.Range("Q6").Select()
